$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the title paragraph's runs so the text reads as one clean
#    sentence and the spell-check "squiggle" markers (w:proofErr) that
#    surrounded each individual word are gone. The character "š" keeps its
#    own run (it carries its own sr-Latn-RS language formatting), the rest
#    of the sentence collapses into two plain runs.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range.Duplicate
$titlePara.Delete()

$titleInsertionPoint = $d.Paragraphs(1).Range.Duplicate
$titleInsertionPoint.Collapse(1)

$titleXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
            "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
            "<pkg:xmlData>" +
            "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
            "<w:body>" +
            "<w:p>" +
            "<w:r><w:t>Definisanje op</w:t></w:r>" +
            "<w:r><w:rPr><w:lang w:val='sr-Latn-RS'/></w:rPr><w:t>&#353;</w:t></w:r>" +
            "<w:r><w:t>tih koraka SCRUM metodologije</w:t></w:r>" +
            "</w:p>" +
            "</w:body>" +
            "</w:document>" +
            "</pkg:xmlData></pkg:part></pkg:package>"

$titleInsertionPoint.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2. Append a blank paragraph followed by a new "Minela test" paragraph at
#    the end of the document.
# ---------------------------------------------------------------------------
$endPoint = $d.Content
$endPoint.Collapse(0)

$tailXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
           "<pkg:xmlData>" +
           "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:body>" +
           "<w:p/>" +
           "<w:p><w:r><w:t>Minela test</w:t></w:r></w:p>" +
           "</w:body>" +
           "</w:document>" +
           "</pkg:xmlData></pkg:part></pkg:package>"

$endPoint.InsertXML($tailXml)
